# Adds an "author" column, normalizes previously-blank/boolean/numeric cells to
# text representations ("N/A", "TRUE"/"FALSE", plain numbers-as-text) so the
# output function no longer breaks on NULL/blank values, and rewrites the
# ids/referenced_works/related_works cells using R-style c()/list() literals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "author" column before the existing "ab" column (column C).
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").Value = "author"

$author = @'
list(au_id = c("https://openalex.org/A5074772971", "https://openalex.org/A5082148123", "https://openalex.org/A5031838322", "https://openalex.org/A5021425074", "https://openalex.org/A5088788713", "https://openalex.org/A5008560207"), au_display_name = c("Ahmed Swidan", "Keith A. Joiner", "Edison Jewson", "Nicolas Carroll", "David Champ", "Gennady Shpak"), au_orcid = c("https://orcid.org/0000-0002-3901-160X", NA, NA, NA, NA, NA), author_position = c("first", "middle", "middle", "middle", "middle", "last"
), au_affiliation_raw = c("Arab Academy for Science, Technology and Maritime Transport, Alexandria, Egypt; University of New South Wales, Canberra, Australia", "University of New South Wales, Canberra, Australia", "University of New South Wales, Canberra, Australia", "University of New South Wales, Canberra, Australia", "University of New South Wales, Canberra, Australia", "University of New South Wales, Canberra, Australia"), institution_id = c("https://openalex.org/I59272784", "https://openalex.org/I188329596", 
"https://openalex.org/I188329596", "https://openalex.org/I188329596", "https://openalex.org/I188329596", "https://openalex.org/I188329596"), institution_display_name = c("Arab Academy for Science, Technology, and Maritime Transport", "University of Canberra", "University of Canberra", "University of Canberra", "University of Canberra", "University of Canberra"), institution_ror = c("https://ror.org/0004vyj87", "https://ror.org/04s1nv328", "https://ror.org/04s1nv328", "https://ror.org/04s1nv328", 
"https://ror.org/04s1nv328", "https://ror.org/04s1nv328"), institution_country_code = c("EG", "AU", "AU", "AU", "AU", "AU"), institution_type = c("education", "education", "education", "education", "education", "education"), institution_lineage = c("https://openalex.org/I59272784", "https://openalex.org/I188329596", "https://openalex.org/I188329596", "https://openalex.org/I188329596", "https://openalex.org/I188329596", "https://openalex.org/I188329596"))
'@
$ws.Range("C2").Value = $author

# Columns that used to be left blank now read "N/A".
$naCols = @("I2","J2","K2","L2","M2","N2","O2","S2","V2")
foreach ($addr in $naCols) {
    $ws.Range($addr).Value = "N/A"
}

# Booleans now round-trip as literal "TRUE"/"FALSE" text instead of native bools.
$boolCols = @("P2","Q2","T2","AE2","AF2")
foreach ($addr in $boolCols) {
    $ws.Range($addr).Value = "'FALSE"
}

# cited_by_count / publication_year now round-trip as text instead of numbers.
$ws.Range("W2").Value = "'1"
$ws.Range("X2").Value = "'2022"

# doi column (previously duplicated the url column) now holds the DOI link.
$ws.Range("AA2").Value = "https://doi.org/10.1109/itc-egypt55520.2022.9855715"

# ids / referenced_works / related_works now serialize as R c()/list() literals.
$ids = @'
c(openalex = "https://openalex.org/W4292348739", doi = "https://doi.org/10.1109/itc-egypt55520.2022.9855715")
'@
$ws.Range("Z2").Value = $ids

$refWorks = @'
c("https://openalex.org/W2023025389", "https://openalex.org/W2025766927", "https://openalex.org/W2068903810", "https://openalex.org/W2072410869", "https://openalex.org/W2616044715", "https://openalex.org/W2625163834", "https://openalex.org/W3091033517", "https://openalex.org/W3142977891")
'@
$ws.Range("AC2").Value = $refWorks

$relWorks = @'
c("https://openalex.org/W2920882006", "https://openalex.org/W1567987063", "https://openalex.org/W2391860589", "https://openalex.org/W2030429945", "https://openalex.org/W2061295771", "https://openalex.org/W4385221012", "https://openalex.org/W4293195837", "https://openalex.org/W114878902", "https://openalex.org/W2936836059", "https://openalex.org/W3092253083")
'@
$ws.Range("AD2").Value = $relWorks

